# Applies refreshed Moogle market-price figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns)
# to the affected leve rows on each profession sheet, as produced
# by the scheduled market-data refresh run. Values only; no
# formulas or formatting are touched.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 32
$ws.Range("H32").Value = 1638.7778
$ws.Range("J32").Value = 1657.5
$ws.Range("L32").Value = 1657.5
$ws.Range("N32").Value = -2309.5

# Row 33
$ws.Range("H33").Value = 36048.8
$ws.Range("I33").Value = 50868.9
$ws.Range("J33").Value = 6408.6
$ws.Range("K33").Value = 50868.9
$ws.Range("L33").Value = 6408.6
$ws.Range("M33").Value = -50639.9
$ws.Range("N33").Value = -6866.6

# Row 125
$ws.Range("H125").Value = 262377.9
$ws.Range("I125").Value = 16950
$ws.Range("J125").Value = 311463.5
$ws.Range("K125").Value = 152550
$ws.Range("L125").Value = 2803171.5
$ws.Range("M125").Value = -150090
$ws.Range("N125").Value = -2808091.5

# Row 126
$ws.Range("H126").Value = 79900
$ws.Range("J126").Value = 79900
$ws.Range("L126").Value = 79900
$ws.Range("N126").Value = -89780

# Row 132
$ws.Range("H132").Value = 3210.093
$ws.Range("I132").Value = 2848.2632
$ws.Range("J132").Value = 5960
$ws.Range("K132").Value = 8544.7896
$ws.Range("L132").Value = 17880
$ws.Range("M132").Value = -6014.7896
$ws.Range("N132").Value = -22940

# Row 138
$ws.Range("H138").Value = 5189.485
$ws.Range("J138").Value = 6627.4375
$ws.Range("L138").Value = 19882.3125
$ws.Range("N138").Value = -30162.3125

# Row 141
$ws.Range("H141").Value = 4466.222
$ws.Range("I141").Value = 3228.8333
$ws.Range("J141").Value = 6941
$ws.Range("K141").Value = 9686.499899999999
$ws.Range("L141").Value = 20823
$ws.Range("M141").Value = -4506.499899999999
$ws.Range("N141").Value = -31183


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 7564.3555
$ws.Range("I32").Value = 2741.5122
$ws.Range("K32").Value = 2741.5122
$ws.Range("M32").Value = -2454.5122

# Row 97
$ws.Range("H97").Value = 1655.2142
$ws.Range("I97").Value = 1984.4546
$ws.Range("J97").Value = 448
$ws.Range("K97").Value = 1984.4546
$ws.Range("L97").Value = 448
$ws.Range("M97").Value = -1488.4546
$ws.Range("N97").Value = -1440

# Row 122
$ws.Range("H122").Value = 2911.2646
$ws.Range("I122").Value = 2954.1724
$ws.Range("K122").Value = 8862.5172
$ws.Range("M122").Value = -6412.5172


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 94
$ws.Range("H94").Value = 973.46155
$ws.Range("I94").Value = 648.3043
$ws.Range("K94").Value = 648.3043
$ws.Range("M94").Value = -197.3043

# Row 99
$ws.Range("H99").Value = 1853.3529
$ws.Range("J99").Value = 2701.6
$ws.Range("L99").Value = 2701.6
$ws.Range("N99").Value = -5697.6

# Row 105
$ws.Range("H105").Value = 2505.4707
$ws.Range("J105").Value = 3550.75
$ws.Range("L105").Value = 3550.75
$ws.Range("N105").Value = -7044.75

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 134
$ws.Range("H134").Value = 2574.96
$ws.Range("J134").Value = 10633.333
$ws.Range("L134").Value = 31899.999
$ws.Range("N134").Value = -36969.999


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 13423.846
$ws.Range("I31").Value = 5401
$ws.Range("J31").Value = 40166.668
$ws.Range("K31").Value = 5401
$ws.Range("L31").Value = 40166.668
$ws.Range("M31").Value = -5106
$ws.Range("N31").Value = -40756.668

# Row 34
$ws.Range("H34").Value = 13423.846
$ws.Range("I34").Value = 5401
$ws.Range("J34").Value = 40166.668
$ws.Range("K34").Value = 5401
$ws.Range("L34").Value = 40166.668
$ws.Range("M34").Value = -5199
$ws.Range("N34").Value = -40570.668

# Row 62
$ws.Range("H62").Value = 2997.25
$ws.Range("I62").Value = 2997.25
$ws.Range("K62").Value = 2997.25
$ws.Range("M62").Value = -2373.25

# Row 65
$ws.Range("H65").Value = 2997.25
$ws.Range("I65").Value = 2997.25
$ws.Range("K65").Value = 14986.25
$ws.Range("M65").Value = -11866.25

# Row 82
$ws.Range("H82").Value = 199950
$ws.Range("J82").Value = 199950
$ws.Range("L82").Value = 199950
$ws.Range("N82").Value = -200672

# Row 85
$ws.Range("H85").Value = 199950
$ws.Range("J85").Value = 199950
$ws.Range("L85").Value = 199950
$ws.Range("N85").Value = -202446

# Row 99
$ws.Range("H99").Value = 1091423.6
$ws.Range("I99").Value = 1503419.4
$ws.Range("K99").Value = 1503419.4
$ws.Range("M99").Value = -1501921.4

# Row 107
$ws.Range("H107").Value = 762
$ws.Range("I107").Value = 619.6111
$ws.Range("K107").Value = 619.6111
$ws.Range("M107").Value = 1300.3889

# Row 126
$ws.Range("H126").Value = 1091423.6
$ws.Range("I126").Value = 1503419.4
$ws.Range("K126").Value = 4510258.199999999
$ws.Range("M126").Value = -4507788.199999999

# Row 134
$ws.Range("H134").Value = 4260.3335
$ws.Range("I134").Value = 2364.9473
$ws.Range("J134").Value = 8761.875
$ws.Range("K134").Value = 7094.841899999999
$ws.Range("L134").Value = 26285.625
$ws.Range("M134").Value = -4559.841899999999
$ws.Range("N134").Value = -31355.625


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 2
$ws.Range("H2").Value = 75.03846
$ws.Range("I2").Value = 88.3125
$ws.Range("J2").Value = 53.8
$ws.Range("K2").Value = 529.875
$ws.Range("L2").Value = 322.8
$ws.Range("M2").Value = -416.875
$ws.Range("N2").Value = -548.8

# Row 12
$ws.Range("H12").Value = 324.16666
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 446.92307
$ws.Range("K12").Value = 15
$ws.Range("L12").Value = 1340.76921
$ws.Range("M12").Value = 158
$ws.Range("N12").Value = -1686.76921

# Row 92
$ws.Range("H92").Value = 88.40000000000001
$ws.Range("I92").Value = 89.75
$ws.Range("J92").Value = 83
$ws.Range("K92").Value = 269.25
$ws.Range("L92").Value = 249
$ws.Range("M92").Value = 978.75
$ws.Range("N92").Value = -2745


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 2
$ws.Range("H2").Value = 694.04
$ws.Range("I2").Value = 953.8333
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 953.8333
$ws.Range("L2").Value = 26
$ws.Range("M2").Value = -840.8333
$ws.Range("N2").Value = -252

# Row 103
$ws.Range("H103").Value = 62821.145
$ws.Range("J103").Value = 62821.145
$ws.Range("L103").Value = 62821.145
$ws.Range("N103").Value = -65165.145

# Row 111
$ws.Range("H111").Value = 199791.67
$ws.Range("J111").Value = 199791.67
$ws.Range("L111").Value = 199791.67
$ws.Range("N111").Value = -205925.67

# Row 122
$ws.Range("H122").Value = 4337.2964
$ws.Range("I122").Value = 3879.7368
$ws.Range("J122").Value = 5424
$ws.Range("K122").Value = 11639.2104
$ws.Range("L122").Value = 16272
$ws.Range("M122").Value = -9189.2104
$ws.Range("N122").Value = -21172

# Row 123
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -34900

# Row 132
$ws.Range("H132").Value = 5312.8647
$ws.Range("I132").Value = 5014.52
$ws.Range("K132").Value = 15043.56
$ws.Range("M132").Value = -12513.56


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 104
$ws.Range("H104").Value = 199950
$ws.Range("J104").Value = 199950
$ws.Range("L104").Value = 199950
$ws.Range("N104").Value = -206938

# Row 132
$ws.Range("H132").Value = 6807.4287
$ws.Range("I132").Value = 6417.3335
$ws.Range("J132").Value = 7327.5557
$ws.Range("K132").Value = 19252.0005
$ws.Range("L132").Value = 21982.6671
$ws.Range("M132").Value = -16722.0005
$ws.Range("N132").Value = -27042.6671

# Row 134
$ws.Range("H134").Value = 68000
$ws.Range("J134").Value = 68000
$ws.Range("L134").Value = 68000
$ws.Range("N134").Value = -78140

# Row 136
$ws.Range("H136").Value = 5478.6665
$ws.Range("I136").Value = 5503.9414
$ws.Range("J136").Value = 5445.615
$ws.Range("K136").Value = 16511.8242
$ws.Range("L136").Value = 16336.845
$ws.Range("M136").Value = -13961.8242
$ws.Range("N136").Value = -21436.845

# Row 137
$ws.Range("H137").Value = 56416
$ws.Range("I137").Value = 20390
$ws.Range("J137").Value = 74429
$ws.Range("K137").Value = 20390
$ws.Range("L137").Value = 74429
$ws.Range("M137").Value = -15290
$ws.Range("N137").Value = -84629
